$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.729.38"
$ws.Range("E2").Value = "  -2.66%  "
$ws.Range("D3").Value = "1.565.40"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'206.30"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("D8").Value = "'21.83"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'0.0584"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("D11").Value = "'0.0863"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "1.787.89"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "1.563.94"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("E14").Value = "  -2.79%  "
$ws.Range("D15").Value = "'0.513"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "26.785.32"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").Value = "'61.32"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").Value = "'214.99"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "'9.31"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").Value = "'152.45"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").Value = "1.392.32"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").Value = "'1.53"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").Value = "'0.927"
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("E38").Value = "  -2.89%  "
$ws.Range("D39").Value = "'0.528"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").Value = "'0.817"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "'0.993"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").Value = "'1.79"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "'63.17"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").Value = "1.700.95"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "'85.73"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").Value = "0.0₇0982"
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("D50").Value = "'0.0951"
$ws.Range("E50").Value = "  -0.85%  "
